$wb = $excel.ActiveWorkbook

# Fix the typo in the shared string used by the Lists sheet (FireFox-Headless -> Firefox-Headless)
$listsWs = $wb.Worksheets.Item("Lists")
$listsWs.Range("A9").Value = "Firefox-Headless"

# Move the selection to Lists!J11 and make Lists the active/selected sheet/tab
$listsWs.Activate()
$listsWs.Range("J11").Select()
